# Fruta / hortaliza, semanal
# Update dates (D) and volumes (M) per row - weekly data reshuffle.
# Rows 2 and 3 also swap their unit-of-sale / price-per-kg / kg-per-unit
# (Q, S, T) values along with the date/volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44327
$ws.Range("M2").Value = 60
$ws.Range("Q2").Value = "$/caja 10 kilos empedrada"
$ws.Range("S2").Value = 11500
$ws.Range("T2").Value = 1

# Row 3
$ws.Range("D3").Value = 44309
$ws.Range("M3").Value = 80
$ws.Range("Q3").Value = "$/caja 14 kilos granel"
$ws.Range("S3").Value = 821
$ws.Range("T3").Value = 14

# Row 4
$ws.Range("D4").Value = 44322
$ws.Range("M4").Value = 60

# Row 6
$ws.Range("D6").Value = 44316
$ws.Range("M6").Value = 120

# Row 7
$ws.Range("D7").Value = 44302
$ws.Range("M7").Value = 80

# Row 8
$ws.Range("D8").Value = 44313
$ws.Range("M8").Value = 120

# Row 9
$ws.Range("D9").Value = 44330
$ws.Range("M9").Value = 60
